$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52: 'Embraer' 'Embraer 190'
$ws.Cells.Item(52, 1).Value = "Embraer"
$ws.Cells.Item(52, 2).Value = "Embraer 190"
$ws.Cells.Item(52, 3).Value = "Regional"
$ws.Cells.Item(52, 4).Value = 2004
$ws.Cells.Item(52, 5).Value = 51843.125
$ws.Cells.Item(52, 6).Value = 44140.625
$ws.Cells.Item(52, 7).Value = 124
$ws.Cells.Item(52, 8).Value = 16347.75
$ws.Cells.Item(52, 9).Value = 16.77016917536357
$ws.Cells.Item(52, 10).Value = 0.3392881252903102
$ws.Cells.Item(52, 11).Value = 1.712372106412916
$ws.Cells.Item(52, 12).Value = 226.7204301075269
$ws.Cells.Item(52, 13).Value = $null
$ws.Cells.Item(52, 14).Value = 9.754665681149657
$ws.Cells.Item(52, 15).Value = 0.04141733060561711
$ws.Cells.Item(52, 16).Value = 30.545
$ws.Cells.Item(52, 17).Value = 0.7117598908588553
$ws.Cells.Item(52, 18).Value = 0.4763098239510975
$ws.Cells.Item(52, 19).Value = $null
$ws.Cells.Item(52, 20).Value = $null
$ws.Cells.Item(52, 21).Value = $null
$ws.Cells.Item(52, 22).Value = $null
$ws.Cells.Item(52, 23).Value = 99
$ws.Cells.Item(52, 24).Value = 10.605
$ws.Cells.Item(52, 25).Value = 7.026875
$ws.Cells.Item(52, 26).Value = 61088.37579710694

# Row 53: 'Embraer' 'Embraer-135'
$ws.Cells.Item(53, 1).Value = "Embraer"
$ws.Cells.Item(53, 2).Value = "Embraer-135"
$ws.Cells.Item(53, 3).Value = "Regional"
$ws.Cells.Item(53, 4).Value = 1999
$ws.Cells.Item(53, 5).Value = 21096
$ws.Cells.Item(53, 6).Value = 15840
$ws.Cells.Item(53, 7).Value = 37
$ws.Cells.Item(53, 8).Value = 8216.8
$ws.Cells.Item(53, 9).Value = 18.16621531661354
$ws.Cells.Item(53, 10).Value = 0.3113135990720869
$ws.Cells.Item(53, 11).Value = $null
$ws.Cells.Item(53, 12).Value = 339.1891891891892
$ws.Cells.Item(53, 13).Value = $null
$ws.Cells.Item(53, 14).Value = 7.846846424384525
$ws.Cells.Item(53, 15).Value = 0.05070665796303604
$ws.Cells.Item(53, 16).Value = 20.04
$ws.Cells.Item(53, 17).Value = 0.5827229105762763
$ws.Cells.Item(53, 18).Value = 0.5342456021720363
$ws.Cells.Item(53, 19).Value = $null
$ws.Cells.Item(53, 20).Value = $null
$ws.Cells.Item(53, 21).Value = $null
$ws.Cells.Item(53, 22).Value = $null
$ws.Cells.Item(53, 23).Value = 37
$ws.Cells.Item(53, 24).Value = 6.76
$ws.Cells.Item(53, 25).Value = 4.762
$ws.Cells.Item(53, 26).Value = 49755.49877845949

# Row 54: 'Embraer' 'Embraer-140'
$ws.Cells.Item(54, 1).Value = "Embraer"
$ws.Cells.Item(54, 2).Value = "Embraer-140"
$ws.Cells.Item(54, 3).Value = "Regional"
$ws.Cells.Item(54, 4).Value = 2001
$ws.Cells.Item(54, 5).Value = $null
$ws.Cells.Item(54, 6).Value = $null
$ws.Cells.Item(54, 7).Value = 44
$ws.Cells.Item(54, 8).Value = $null
$ws.Cells.Item(54, 9).Value = $null
$ws.Cells.Item(54, 10).Value = $null
$ws.Cells.Item(54, 11).Value = $null
$ws.Cells.Item(54, 12).Value = 268.1818181818182
$ws.Cells.Item(54, 13).Value = $null
$ws.Cells.Item(54, 14).Value = $null
$ws.Cells.Item(54, 15).Value = $null
$ws.Cells.Item(54, 16).Value = $null
$ws.Cells.Item(54, 17).Value = $null
$ws.Cells.Item(54, 18).Value = $null
$ws.Cells.Item(54, 19).Value = $null
$ws.Cells.Item(54, 20).Value = $null
$ws.Cells.Item(54, 21).Value = $null
$ws.Cells.Item(54, 22).Value = $null
$ws.Cells.Item(54, 23).Value = $null
$ws.Cells.Item(54, 24).Value = $null
$ws.Cells.Item(54, 25).Value = $null
$ws.Cells.Item(54, 26).Value = $null

# Row 55: 'Embraer ' 'EMB-120 Brasilia'
$ws.Cells.Item(55, 1).Value = "Embraer "
$ws.Cells.Item(55, 2).Value = "EMB-120 Brasilia"
$ws.Cells.Item(55, 3).Value = "Regional"
$ws.Cells.Item(55, 4).Value = 1985
$ws.Cells.Item(55, 5).Value = $null
$ws.Cells.Item(55, 6).Value = $null
$ws.Cells.Item(55, 7).Value = 30
$ws.Cells.Item(55, 8).Value = $null
$ws.Cells.Item(55, 9).Value = $null
$ws.Cells.Item(55, 10).Value = $null
$ws.Cells.Item(55, 11).Value = $null
$ws.Cells.Item(55, 12).Value = 252
$ws.Cells.Item(55, 13).Value = $null
$ws.Cells.Item(55, 14).Value = $null
$ws.Cells.Item(55, 15).Value = $null
$ws.Cells.Item(55, 16).Value = $null
$ws.Cells.Item(55, 17).Value = $null
$ws.Cells.Item(55, 18).Value = $null
$ws.Cells.Item(55, 19).Value = $null
$ws.Cells.Item(55, 20).Value = $null
$ws.Cells.Item(55, 21).Value = $null
$ws.Cells.Item(55, 22).Value = $null
$ws.Cells.Item(55, 23).Value = $null
$ws.Cells.Item(55, 24).Value = $null
$ws.Cells.Item(55, 25).Value = $null
$ws.Cells.Item(55, 26).Value = $null

# Row 56: 'Embraer ' 'Embraer ERJ-175'
$ws.Cells.Item(56, 1).Value = "Embraer "
$ws.Cells.Item(56, 2).Value = "Embraer ERJ-175"
$ws.Cells.Item(56, 3).Value = "Regional"
$ws.Cells.Item(56, 4).Value = 2005
$ws.Cells.Item(56, 5).Value = 37500
$ws.Cells.Item(56, 6).Value = 31700
$ws.Cells.Item(56, 7).Value = 88
$ws.Cells.Item(56, 8).Value = 11625
$ws.Cells.Item(56, 9).Value = 18.19340342917936
$ws.Cells.Item(56, 10).Value = 0.3108455349745253
$ws.Cells.Item(56, 11).Value = 1.528589461364119
$ws.Cells.Item(56, 12).Value = 247.8409090909091
$ws.Cells.Item(56, 13).Value = $null
$ws.Cells.Item(56, 14).Value = 11.28743811881188
$ws.Cells.Item(56, 15).Value = 0.03525045750342684
$ws.Cells.Item(56, 16).Value = 28.65
$ws.Cells.Item(56, 17).Value = 0.7270491548603031
$ws.Cells.Item(56, 18).Value = 0.4275455547427663
$ws.Cells.Item(56, 19).Value = $null
$ws.Cells.Item(56, 20).Value = $null
$ws.Cells.Item(56, 21).Value = $null
$ws.Cells.Item(56, 22).Value = $null
$ws.Cells.Item(56, 23).Value = 75
$ws.Cells.Item(56, 24).Value = 9.82
$ws.Cells.Item(56, 25).Value = 5.13
$ws.Cells.Item(56, 26).Value = 37180.61653102213

# Row 57: 'Embraer ' 'Embraer-145'
$ws.Cells.Item(57, 1).Value = "Embraer "
$ws.Cells.Item(57, 2).Value = "Embraer-145"
$ws.Cells.Item(57, 3).Value = "Regional"
$ws.Cells.Item(57, 4).Value = 1996
$ws.Cells.Item(57, 5).Value = 20016.66666666667
$ws.Cells.Item(57, 6).Value = 17333.33333333333
$ws.Cells.Item(57, 7).Value = 50
$ws.Cells.Item(57, 8).Value = 5528
$ws.Cells.Item(57, 9).Value = 18.22195069609804
$ws.Cells.Item(57, 10).Value = 0.3103601476304506
$ws.Cells.Item(57, 11).Value = 2.103581856805667
$ws.Cells.Item(57, 12).Value = 244.1333333333333
$ws.Cells.Item(57, 13).Value = $null
$ws.Cells.Item(57, 14).Value = 7.975146541617818
$ws.Cells.Item(57, 15).Value = 0.04995164617426398
$ws.Cells.Item(57, 16).Value = 20.2
$ws.Cells.Item(57, 17).Value = 0.6081109382651159
$ws.Cells.Item(57, 18).Value = 0.5103718988149327
$ws.Cells.Item(57, 19).Value = $null
$ws.Cells.Item(57, 20).Value = $null
$ws.Cells.Item(57, 21).Value = $null
$ws.Cells.Item(57, 22).Value = $null
$ws.Cells.Item(57, 23).Value = 50
$ws.Cells.Item(57, 24).Value = 6.75
$ws.Cells.Item(57, 25).Value = 4.837777777777777
$ws.Cells.Item(57, 26).Value = 44777.31705387121

# Row 58: 'Fokker ' 'Fokker 100'
$ws.Cells.Item(58, 1).Value = "Fokker "
$ws.Cells.Item(58, 2).Value = "Fokker 100"
$ws.Cells.Item(58, 3).Value = "Regional"
$ws.Cells.Item(58, 4).Value = 1988
$ws.Cells.Item(58, 5).Value = $null
$ws.Cells.Item(58, 6).Value = $null
$ws.Cells.Item(58, 7).Value = 122
$ws.Cells.Item(58, 8).Value = $null
$ws.Cells.Item(58, 9).Value = $null
$ws.Cells.Item(58, 10).Value = $null
$ws.Cells.Item(58, 11).Value = $null
$ws.Cells.Item(58, 12).Value = 202.3142076502733
$ws.Cells.Item(58, 13).Value = $null
$ws.Cells.Item(58, 14).Value = $null
$ws.Cells.Item(58, 15).Value = $null
$ws.Cells.Item(58, 16).Value = $null
$ws.Cells.Item(58, 17).Value = $null
$ws.Cells.Item(58, 18).Value = $null
$ws.Cells.Item(58, 19).Value = $null
$ws.Cells.Item(58, 20).Value = $null
$ws.Cells.Item(58, 21).Value = $null
$ws.Cells.Item(58, 22).Value = $null
$ws.Cells.Item(58, 23).Value = $null
$ws.Cells.Item(58, 24).Value = $null
$ws.Cells.Item(58, 25).Value = $null
$ws.Cells.Item(58, 26).Value = $null

# Row 59: 'Gates Learjet ' 'Lear-31/35/36'
$ws.Cells.Item(59, 1).Value = "Gates Learjet "
$ws.Cells.Item(59, 2).Value = "Lear-31/35/36"
$ws.Cells.Item(59, 3).Value = "Regional"
$ws.Cells.Item(59, 4).Value = 1974
$ws.Cells.Item(59, 5).Value = $null
$ws.Cells.Item(59, 6).Value = $null
$ws.Cells.Item(59, 7).Value = $null
$ws.Cells.Item(59, 8).Value = $null
$ws.Cells.Item(59, 9).Value = $null
$ws.Cells.Item(59, 10).Value = $null
$ws.Cells.Item(59, 11).Value = $null
$ws.Cells.Item(59, 12).Value = $null
$ws.Cells.Item(59, 13).Value = $null
$ws.Cells.Item(59, 14).Value = $null
$ws.Cells.Item(59, 15).Value = $null
$ws.Cells.Item(59, 16).Value = $null
$ws.Cells.Item(59, 17).Value = $null
$ws.Cells.Item(59, 18).Value = $null
$ws.Cells.Item(59, 19).Value = $null
$ws.Cells.Item(59, 20).Value = $null
$ws.Cells.Item(59, 21).Value = $null
$ws.Cells.Item(59, 22).Value = $null
$ws.Cells.Item(59, 23).Value = $null
$ws.Cells.Item(59, 24).Value = $null
$ws.Cells.Item(59, 25).Value = $null
$ws.Cells.Item(59, 26).Value = $null

# Row 60: 'Lockheed' 'L1011-1/100/200'
$ws.Cells.Item(60, 1).Value = "Lockheed"
$ws.Cells.Item(60, 2).Value = "L1011-1/100/200"
$ws.Cells.Item(60, 3).Value = "Wide"
$ws.Cells.Item(60, 4).Value = 1973
$ws.Cells.Item(60, 5).Value = 224982
$ws.Cells.Item(60, 6).Value = 153314
$ws.Cells.Item(60, 7).Value = 400
$ws.Cells.Item(60, 8).Value = 120705
$ws.Cells.Item(60, 9).Value = 17.5682
$ws.Cells.Item(60, 10).Value = 0.3219043875967976
$ws.Cells.Item(60, 11).Value = 1.69015679246197
$ws.Cells.Item(60, 12).Value = 277.7625
$ws.Cells.Item(60, 13).Value = 14.3153
$ws.Cells.Item(60, 14).Value = $null
$ws.Cells.Item(60, 15).Value = $null
$ws.Cells.Item(60, 16).Value = $null
$ws.Cells.Item(60, 17).Value = 0.7611246338515383
$ws.Cells.Item(60, 18).Value = 0.4229325570082481
$ws.Cells.Item(60, 19).Value = $null
$ws.Cells.Item(60, 20).Value = $null
$ws.Cells.Item(60, 21).Value = $null
$ws.Cells.Item(60, 22).Value = $null
$ws.Cells.Item(60, 23).Value = 314
$ws.Cells.Item(60, 24).Value = $null
$ws.Cells.Item(60, 25).Value = 4.54
$ws.Cells.Item(60, 26).Value = 163349.5957069433

# Row 61: 'Lockheed ' 'L1011-500'
$ws.Cells.Item(61, 1).Value = "Lockheed "
$ws.Cells.Item(61, 2).Value = "L1011-500"
$ws.Cells.Item(61, 3).Value = "Wide"
$ws.Cells.Item(61, 4).Value = 1979
$ws.Cells.Item(61, 5).Value = $null
$ws.Cells.Item(61, 6).Value = $null
$ws.Cells.Item(61, 7).Value = 330
$ws.Cells.Item(61, 8).Value = $null
$ws.Cells.Item(61, 9).Value = $null
$ws.Cells.Item(61, 10).Value = $null
$ws.Cells.Item(61, 11).Value = 1.8209
$ws.Cells.Item(61, 12).Value = 330.0757575757576
$ws.Cells.Item(61, 13).Value = $null
$ws.Cells.Item(61, 14).Value = $null
$ws.Cells.Item(61, 15).Value = $null
$ws.Cells.Item(61, 16).Value = $null
$ws.Cells.Item(61, 17).Value = $null
$ws.Cells.Item(61, 18).Value = $null
$ws.Cells.Item(61, 19).Value = $null
$ws.Cells.Item(61, 20).Value = $null
$ws.Cells.Item(61, 21).Value = $null
$ws.Cells.Item(61, 22).Value = $null
$ws.Cells.Item(61, 23).Value = $null
$ws.Cells.Item(61, 24).Value = $null
$ws.Cells.Item(61, 25).Value = $null
$ws.Cells.Item(61, 26).Value = $null

# Row 62: 'McDonnell Douglas' 'DC10-30'
$ws.Cells.Item(62, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(62, 2).Value = "DC10-30"
$ws.Cells.Item(62, 3).Value = "Wide"
$ws.Cells.Item(62, 4).Value = 1972
$ws.Cells.Item(62, 5).Value = $null
$ws.Cells.Item(62, 6).Value = $null
$ws.Cells.Item(62, 7).Value = 380
$ws.Cells.Item(62, 8).Value = $null
$ws.Cells.Item(62, 9).Value = 17.7708
$ws.Cells.Item(62, 10).Value = 0.3182344442668905
$ws.Cells.Item(62, 11).Value = 1.803152741832779
$ws.Cells.Item(62, 12).Value = 318.9342105263158
$ws.Cells.Item(62, 13).Value = 15.2348
$ws.Cells.Item(62, 14).Value = $null
$ws.Cells.Item(62, 15).Value = $null
$ws.Cells.Item(62, 16).Value = $null
$ws.Cells.Item(62, 17).Value = $null
$ws.Cells.Item(62, 18).Value = $null
$ws.Cells.Item(62, 19).Value = $null
$ws.Cells.Item(62, 20).Value = $null
$ws.Cells.Item(62, 21).Value = $null
$ws.Cells.Item(62, 22).Value = $null
$ws.Cells.Item(62, 23).Value = 266
$ws.Cells.Item(62, 24).Value = $null
$ws.Cells.Item(62, 25).Value = $null
$ws.Cells.Item(62, 26).Value = 176259.1642868837

# Row 63: 'McDonnell Douglas' 'DC10-40'
$ws.Cells.Item(63, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(63, 2).Value = "DC10-40"
$ws.Cells.Item(63, 3).Value = "Wide"
$ws.Cells.Item(63, 4).Value = 1972
$ws.Cells.Item(63, 5).Value = 256280
$ws.Cells.Item(63, 6).Value = 177355
$ws.Cells.Item(63, 7).Value = 380
$ws.Cells.Item(63, 8).Value = 137520
$ws.Cells.Item(63, 9).Value = 17.4331
$ws.Cells.Item(63, 10).Value = 0.3243990261157257
$ws.Cells.Item(63, 11).Value = 1.815799269772046
$ws.Cells.Item(63, 12).Value = 323.5526315789473
$ws.Cells.Item(63, 13).Value = 13.9129
$ws.Cells.Item(63, 14).Value = $null
$ws.Cells.Item(63, 15).Value = $null
$ws.Cells.Item(63, 16).Value = $null
$ws.Cells.Item(63, 17).Value = 0.7498489997700178
$ws.Cells.Item(63, 18).Value = 0.4326919339257812
$ws.Cells.Item(63, 19).Value = $null
$ws.Cells.Item(63, 20).Value = $null
$ws.Cells.Item(63, 21).Value = $null
$ws.Cells.Item(63, 22).Value = $null
$ws.Cells.Item(63, 23).Value = 292
$ws.Cells.Item(63, 24).Value = $null
$ws.Cells.Item(63, 25).Value = 4.961538461538462
$ws.Cells.Item(63, 26).Value = 176295.7276686124

# Row 64: 'McDonnell Douglas' 'DC9-30'
$ws.Cells.Item(64, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(64, 2).Value = "DC9-30"
$ws.Cells.Item(64, 3).Value = "Narrow"
$ws.Cells.Item(64, 4).Value = 1966
$ws.Cells.Item(64, 5).Value = 47627
$ws.Cells.Item(64, 6).Value = 39463
$ws.Cells.Item(64, 7).Value = 127
$ws.Cells.Item(64, 8).Value = 13926
$ws.Cells.Item(64, 9).Value = 22.4331
$ws.Cells.Item(64, 10).Value = 0.2520953707770241
$ws.Cells.Item(64, 11).Value = 2.396593879999583
$ws.Cells.Item(64, 12).Value = 203.0629921259843
$ws.Cells.Item(64, 13).Value = 13.798
$ws.Cells.Item(64, 14).Value = $null
$ws.Cells.Item(64, 15).Value = $null
$ws.Cells.Item(64, 16).Value = 28.44
$ws.Cells.Item(64, 17).Value = $null
$ws.Cells.Item(64, 18).Value = $null
$ws.Cells.Item(64, 19).Value = $null
$ws.Cells.Item(64, 20).Value = $null
$ws.Cells.Item(64, 21).Value = $null
$ws.Cells.Item(64, 22).Value = $null
$ws.Cells.Item(64, 23).Value = 100
$ws.Cells.Item(64, 24).Value = 8.4
$ws.Cells.Item(64, 25).Value = 1.038181818181818
$ws.Cells.Item(64, 26).Value = 50239.33056833609

# Row 65: 'McDonnell Douglas' 'DC9-40'
$ws.Cells.Item(65, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(65, 2).Value = "DC9-40"
$ws.Cells.Item(65, 3).Value = "Narrow"
$ws.Cells.Item(65, 4).Value = 1968
$ws.Cells.Item(65, 5).Value = 51710
$ws.Cells.Item(65, 6).Value = 42184
$ws.Cells.Item(65, 7).Value = 128
$ws.Cells.Item(65, 8).Value = 13926
$ws.Cells.Item(65, 9).Value = 22.98764963035378
$ws.Cells.Item(65, 10).Value = 0.2460964991804533
$ws.Cells.Item(65, 11).Value = 2.0777
$ws.Cells.Item(65, 12).Value = 217.3515625
$ws.Cells.Item(65, 13).Value = $null
$ws.Cells.Item(65, 14).Value = $null
$ws.Cells.Item(65, 15).Value = $null
$ws.Cells.Item(65, 16).Value = 28.44
$ws.Cells.Item(65, 17).Value = $null
$ws.Cells.Item(65, 18).Value = $null
$ws.Cells.Item(65, 19).Value = $null
$ws.Cells.Item(65, 20).Value = $null
$ws.Cells.Item(65, 21).Value = $null
$ws.Cells.Item(65, 22).Value = $null
$ws.Cells.Item(65, 23).Value = $null
$ws.Cells.Item(65, 24).Value = 8.5
$ws.Cells.Item(65, 25).Value = 1.036428571428571
$ws.Cells.Item(65, 26).Value = $null

# Row 66: 'McDonnell Douglas' 'DC9-50'
$ws.Cells.Item(66, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(66, 2).Value = "DC9-50"
$ws.Cells.Item(66, 3).Value = "Narrow"
$ws.Cells.Item(66, 4).Value = 1976
$ws.Cells.Item(66, 5).Value = 54885
$ws.Cells.Item(66, 6).Value = 44679
$ws.Cells.Item(66, 7).Value = 139
$ws.Cells.Item(66, 8).Value = 13926
$ws.Cells.Item(66, 9).Value = 23.15547775838727
$ws.Cells.Item(66, 10).Value = 0.2443182490108967
$ws.Cells.Item(66, 11).Value = 2.1014
$ws.Cells.Item(66, 12).Value = 211.0503597122302
$ws.Cells.Item(66, 13).Value = $null
$ws.Cells.Item(66, 14).Value = $null
$ws.Cells.Item(66, 15).Value = $null
$ws.Cells.Item(66, 16).Value = 28.45
$ws.Cells.Item(66, 17).Value = $null
$ws.Cells.Item(66, 18).Value = $null
$ws.Cells.Item(66, 19).Value = $null
$ws.Cells.Item(66, 20).Value = $null
$ws.Cells.Item(66, 21).Value = $null
$ws.Cells.Item(66, 22).Value = $null
$ws.Cells.Item(66, 23).Value = $null
$ws.Cells.Item(66, 24).Value = 8.6
$ws.Cells.Item(66, 25).Value = 1.04375
$ws.Cells.Item(66, 26).Value = $null

# Row 67: 'McDonnell Douglas' 'MD-90'
$ws.Cells.Item(67, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(67, 2).Value = "MD-90"
$ws.Cells.Item(67, 3).Value = "Narrow"
$ws.Cells.Item(67, 4).Value = 1995
$ws.Cells.Item(67, 5).Value = 70760
$ws.Cells.Item(67, 6).Value = 58967
$ws.Cells.Item(67, 7).Value = 167
$ws.Cells.Item(67, 8).Value = 22104
$ws.Cells.Item(67, 9).Value = 16.99770822716683
$ws.Cells.Item(67, 10).Value = 0.3327206303053079
$ws.Cells.Item(67, 11).Value = 1.340189934684793
$ws.Cells.Item(67, 12).Value = 242.5489021956088
$ws.Cells.Item(67, 13).Value = 12.54450404243209
$ws.Cells.Item(67, 14).Value = $null
$ws.Cells.Item(67, 15).Value = $null
$ws.Cells.Item(67, 16).Value = 32.87
$ws.Cells.Item(67, 17).Value = 0.7458748318947507
$ws.Cells.Item(67, 18).Value = 0.4460862698857084
$ws.Cells.Item(67, 19).Value = $null
$ws.Cells.Item(67, 20).Value = $null
$ws.Cells.Item(67, 21).Value = $null
$ws.Cells.Item(67, 22).Value = $null
$ws.Cells.Item(67, 23).Value = 148
$ws.Cells.Item(67, 24).Value = 9.3
$ws.Cells.Item(67, 25).Value = 4.74
$ws.Cells.Item(67, 26).Value = 63894.06446530207

# Row 68: 'McDonnell Douglas' 'MD80/DC9-80'
$ws.Cells.Item(68, 1).Value = "McDonnell Douglas"
$ws.Cells.Item(68, 2).Value = "MD80/DC9-80"
$ws.Cells.Item(68, 3).Value = "Narrow"
$ws.Cells.Item(68, 4).Value = 1980
$ws.Cells.Item(68, 5).Value = 66170.71428571429
$ws.Cells.Item(68, 6).Value = 54258.33333333334
$ws.Cells.Item(68, 7).Value = 172
$ws.Cells.Item(68, 8).Value = 22128.57142857143
$ws.Cells.Item(68, 9).Value = 20.744
$ws.Cells.Item(68, 10).Value = 0.2726224769657761
$ws.Cells.Item(68, 11).Value = 1.654110373626076
$ws.Cells.Item(68, 12).Value = 201.4290697674419
$ws.Cells.Item(68, 13).Value = 13.9129
$ws.Cells.Item(68, 14).Value = 9.609283170080145
$ws.Cells.Item(68, 15).Value = 0.04140655974928667
$ws.Cells.Item(68, 16).Value = 32.85
$ws.Cells.Item(68, 17).Value = $null
$ws.Cells.Item(68, 18).Value = $null
$ws.Cells.Item(68, 19).Value = 0.4818759914491247
$ws.Cells.Item(68, 20).Value = 0.04330909613574976
$ws.Cells.Item(68, 21).Value = 0.009614788308106073
$ws.Cells.Item(68, 22).Value = 0.03369430782764369
$ws.Cells.Item(68, 23).Value = 144
$ws.Cells.Item(68, 24).Value = 9
$ws.Cells.Item(68, 25).Value = 1.725714285714286
$ws.Cells.Item(68, 26).Value = 58284.64332602391

# Row 69: 'McDonnell Douglas ' 'DC10-10'
$ws.Cells.Item(69, 1).Value = "McDonnell Douglas "
$ws.Cells.Item(69, 2).Value = "DC10-10"
$ws.Cells.Item(69, 3).Value = "Wide"
$ws.Cells.Item(69, 4).Value = 1970
$ws.Cells.Item(69, 5).Value = $null
$ws.Cells.Item(69, 6).Value = $null
$ws.Cells.Item(69, 7).Value = 380
$ws.Cells.Item(69, 8).Value = $null
$ws.Cells.Item(69, 9).Value = 17.0953
$ws.Cells.Item(69, 10).Value = 0.3308090915150982
$ws.Cells.Item(69, 11).Value = 1.665838405572881
$ws.Cells.Item(69, 12).Value = 286.6842105263158
$ws.Cells.Item(69, 13).Value = 14.2003
$ws.Cells.Item(69, 14).Value = $null
$ws.Cells.Item(69, 15).Value = $null
$ws.Cells.Item(69, 16).Value = $null
$ws.Cells.Item(69, 17).Value = $null
$ws.Cells.Item(69, 18).Value = $null
$ws.Cells.Item(69, 19).Value = $null
$ws.Cells.Item(69, 20).Value = $null
$ws.Cells.Item(69, 21).Value = $null
$ws.Cells.Item(69, 22).Value = $null
$ws.Cells.Item(69, 23).Value = 284
$ws.Cells.Item(69, 24).Value = $null
$ws.Cells.Item(69, 25).Value = $null
$ws.Cells.Item(69, 26).Value = 169567.8746222493

# Row 70: 'McDonnell Douglas ' 'DC9-10'
$ws.Cells.Item(70, 1).Value = "McDonnell Douglas "
$ws.Cells.Item(70, 2).Value = "DC9-10"
$ws.Cells.Item(70, 3).Value = "Narrow"
$ws.Cells.Item(70, 4).Value = 1965
$ws.Cells.Item(70, 5).Value = 41141
$ws.Cells.Item(70, 6).Value = 33566
$ws.Cells.Item(70, 7).Value = 109
$ws.Cells.Item(70, 8).Value = 13979
$ws.Cells.Item(70, 9).Value = 22.88492472710245
$ws.Cells.Item(70, 10).Value = 0.2472409829726473
$ws.Cells.Item(70, 11).Value = 2.9899
$ws.Cells.Item(70, 12).Value = 204.5871559633028
$ws.Cells.Item(70, 13).Value = $null
$ws.Cells.Item(70, 14).Value = $null
$ws.Cells.Item(70, 15).Value = $null
$ws.Cells.Item(70, 16).Value = 27.25
$ws.Cells.Item(70, 17).Value = $null
$ws.Cells.Item(70, 18).Value = $null
$ws.Cells.Item(70, 19).Value = $null
$ws.Cells.Item(70, 20).Value = $null
$ws.Cells.Item(70, 21).Value = $null
$ws.Cells.Item(70, 22).Value = $null
$ws.Cells.Item(70, 23).Value = $null
$ws.Cells.Item(70, 24).Value = 8.4
$ws.Cells.Item(70, 25).Value = 1.032222222222222
$ws.Cells.Item(70, 26).Value = $null

# Row 71: 'McDonnell Douglas ' 'MD-11'
$ws.Cells.Item(71, 1).Value = "McDonnell Douglas "
$ws.Cells.Item(71, 2).Value = "MD-11"
$ws.Cells.Item(71, 3).Value = "Wide"
$ws.Cells.Item(71, 4).Value = 1990
$ws.Cells.Item(71, 5).Value = 276691
$ws.Cells.Item(71, 6).Value = 195045
$ws.Cells.Item(71, 7).Value = 410
$ws.Cells.Item(71, 8).Value = 144782
$ws.Cells.Item(71, 9).Value = 17.17456583460847
$ws.Cells.Item(71, 10).Value = 0.329294025279972
$ws.Cells.Item(71, 11).Value = 1.744608638314735
$ws.Cells.Item(71, 12).Value = 317.4756097560976
$ws.Cells.Item(71, 13).Value = 17.18295849059495
$ws.Cells.Item(71, 14).Value = $null
$ws.Cells.Item(71, 15).Value = $null
$ws.Cells.Item(71, 16).Value = $null
$ws.Cells.Item(71, 17).Value = 0.7918255422774546
$ws.Cells.Item(71, 18).Value = 0.4158827258454688
$ws.Cells.Item(71, 19).Value = $null
$ws.Cells.Item(71, 20).Value = $null
$ws.Cells.Item(71, 21).Value = $null
$ws.Cells.Item(71, 22).Value = $null
$ws.Cells.Item(71, 23).Value = 254
$ws.Cells.Item(71, 24).Value = $null
$ws.Cells.Item(71, 25).Value = 4.8
$ws.Cells.Item(71, 26).Value = 160911.0912490309

# Row 72: 'Saab-Fairchild ' '340/B'
$ws.Cells.Item(72, 1).Value = "Saab-Fairchild "
$ws.Cells.Item(72, 2).Value = "340/B"
$ws.Cells.Item(72, 3).Value = "Regional"
$ws.Cells.Item(72, 4).Value = 1984
$ws.Cells.Item(72, 5).Value = $null
$ws.Cells.Item(72, 6).Value = $null
$ws.Cells.Item(72, 7).Value = 34
$ws.Cells.Item(72, 8).Value = $null
$ws.Cells.Item(72, 9).Value = $null
$ws.Cells.Item(72, 10).Value = $null
$ws.Cells.Item(72, 11).Value = $null
$ws.Cells.Item(72, 12).Value = 241.9117647058823
$ws.Cells.Item(72, 13).Value = $null
$ws.Cells.Item(72, 14).Value = $null
$ws.Cells.Item(72, 15).Value = $null
$ws.Cells.Item(72, 16).Value = $null
$ws.Cells.Item(72, 17).Value = $null
$ws.Cells.Item(72, 18).Value = $null
$ws.Cells.Item(72, 19).Value = $null
$ws.Cells.Item(72, 20).Value = $null
$ws.Cells.Item(72, 21).Value = $null
$ws.Cells.Item(72, 22).Value = $null
$ws.Cells.Item(72, 23).Value = $null
$ws.Cells.Item(72, 24).Value = $null
$ws.Cells.Item(72, 25).Value = $null
$ws.Cells.Item(72, 26).Value = $null

# Row 73: 'de Havilland' 'Comet 1'
$ws.Cells.Item(73, 1).Value = "de Havilland"
$ws.Cells.Item(73, 2).Value = "Comet 1"
$ws.Cells.Item(73, 3).Value = "Narrow"
$ws.Cells.Item(73, 4).Value = 1952
$ws.Cells.Item(73, 5).Value = $null
$ws.Cells.Item(73, 6).Value = $null
$ws.Cells.Item(73, 7).Value = 44
$ws.Cells.Item(73, 8).Value = $null
$ws.Cells.Item(73, 9).Value = 28.9
$ws.Cells.Item(73, 10).Value = 0.1956844519784796
$ws.Cells.Item(73, 11).Value = 8.624957727272728
$ws.Cells.Item(73, 12).Value = 721.5909090909091
$ws.Cells.Item(73, 13).Value = 13.03220058392828
$ws.Cells.Item(73, 14).Value = $null
$ws.Cells.Item(73, 15).Value = $null
$ws.Cells.Item(73, 16).Value = $null
$ws.Cells.Item(73, 17).Value = $null
$ws.Cells.Item(73, 18).Value = $null
$ws.Cells.Item(73, 19).Value = $null
$ws.Cells.Item(73, 20).Value = $null
$ws.Cells.Item(73, 21).Value = $null
$ws.Cells.Item(73, 22).Value = $null
$ws.Cells.Item(73, 23).Value = $null
$ws.Cells.Item(73, 24).Value = $null
$ws.Cells.Item(73, 25).Value = $null
$ws.Cells.Item(73, 26).Value = $null

# Row 74: 'de Havilland' 'Comet 4'
$ws.Cells.Item(74, 1).Value = "de Havilland"
$ws.Cells.Item(74, 2).Value = "Comet 4"
$ws.Cells.Item(74, 3).Value = "Narrow"
$ws.Cells.Item(74, 4).Value = 1958
$ws.Cells.Item(74, 5).Value = $null
$ws.Cells.Item(74, 6).Value = $null
$ws.Cells.Item(74, 7).Value = 109
$ws.Cells.Item(74, 8).Value = $null
$ws.Cells.Item(74, 9).Value = 26.4
$ws.Cells.Item(74, 10).Value = 0.2142151765976538
$ws.Cells.Item(74, 11).Value = 3.451363636363636
$ws.Cells.Item(74, 12).Value = 313.8715596330275
$ws.Cells.Item(74, 13).Value = 12.94289328120862
$ws.Cells.Item(74, 14).Value = $null
$ws.Cells.Item(74, 15).Value = $null
$ws.Cells.Item(74, 16).Value = $null
$ws.Cells.Item(74, 17).Value = $null
$ws.Cells.Item(74, 18).Value = $null
$ws.Cells.Item(74, 19).Value = $null
$ws.Cells.Item(74, 20).Value = $null
$ws.Cells.Item(74, 21).Value = $null
$ws.Cells.Item(74, 22).Value = $null
$ws.Cells.Item(74, 23).Value = $null
$ws.Cells.Item(74, 24).Value = $null
$ws.Cells.Item(74, 25).Value = $null
$ws.Cells.Item(74, 26).Value = $null

Write-Host "done"
